$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 97) with the August/September SSA data point
$row = 97

$dateCell = $ws.Cells.Item($row, 1)
$dateCell.Formula = "=""2020-09-04"""
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Cells.Item($row, 2).Value = 623090
$ws.Cells.Item($row, 3).Value = 700502
$ws.Cells.Item($row, 4).Value = 85792
$ws.Cells.Item($row, 5).Value = 66851
$ws.Cells.Item($row, 6).Value = 25.17
